$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.233.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.50%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.910.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.42%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +0.22%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.7392"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.91%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'244.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.83%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.23%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3134"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.67%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'27.03"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -4.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.06967"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.13%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.7804"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07976"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.57%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.938.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.290"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.90%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'91.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -3.50%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'14.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.55%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.256.36"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.45%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'5.932"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.01%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'244.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -4.79%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.000007840"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -2.57%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.168.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.79%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.17%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.25%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.655"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -2.53%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.445"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.00%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'165.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +0.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.08%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.1270"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.72%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.108"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -9.16%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.355"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.26%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  +0.99%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.319"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.05%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.083"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.05184"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.90%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.296"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.18%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7464"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'2.763"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.25%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01941"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -1.80%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'2.794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.70%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.355"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.04%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'75.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.90%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.4485"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.17%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'1.940"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.05%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  -0.32%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'7.705"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'9.895"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.24%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -1.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'37.36"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.1201"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.63%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'939.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.75%  "
$ws.Range("E51").Style = "Normal"
